$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 42687508
$ws.Range("I62").Value = 17862122
$ws.Range("K62").Value = 17862122
$ws.Range("M62").Value = -17861498

$ws.Range("H64").Value = 3242.7778
$ws.Range("I64").Value = 3397
$ws.Range("J64").Value = 3050
$ws.Range("K64").Value = 3397
$ws.Range("L64").Value = 3050
$ws.Range("M64").Value = -3149
$ws.Range("N64").Value = -3546

$ws.Range("H65").Value = 42687508
$ws.Range("I65").Value = 17862122
$ws.Range("K65").Value = 89310610
$ws.Range("M65").Value = -89307490

$ws.Range("H67").Value = 3242.7778
$ws.Range("I67").Value = 3397
$ws.Range("J67").Value = 3050
$ws.Range("K67").Value = 3397
$ws.Range("L67").Value = 3050
$ws.Range("M67").Value = -2539
$ws.Range("N67").Value = -4766

$ws.Range("H76").Value = 125003070
$ws.Range("I76").Value = 500001540
$ws.Range("J76").Value = 3583.3333
$ws.Range("K76").Value = 500001540
$ws.Range("L76").Value = 3583.3333
$ws.Range("M76").Value = -500001225
$ws.Range("N76").Value = -4213.3333

$ws.Range("H79").Value = 125003070
$ws.Range("I79").Value = 500001540
$ws.Range("J79").Value = 3583.3333
$ws.Range("K79").Value = 500001540
$ws.Range("L79").Value = 3583.3333
$ws.Range("M79").Value = -500000448
$ws.Range("N79").Value = -5767.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2094432.2
$ws.Range("I61").Value = 1042577.06
$ws.Range("J61").Value = 7353708
$ws.Range("K61").Value = 1042577.06
$ws.Range("L61").Value = 7353708
$ws.Range("M61").Value = -1042365.06
$ws.Range("N61").Value = -7354132

$ws.Range("H88").Value = 9625
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9625
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = -10437

$ws.Range("H91").Value = 9625
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9625
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = -12433

$ws.Range("H97").Value = 438.9
$ws.Range("I97").Value = 454.14285
$ws.Range("J97").Value = 403.33334
$ws.Range("K97").Value = 454.14285
$ws.Range("L97").Value = 403.33334
$ws.Range("M97").Value = 41.85714999999999
$ws.Range("N97").Value = -1395.33334

$ws.Range("H122").Value = 2293.3333
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 3368
$ws.Range("K122").Value = 2850
$ws.Range("L122").Value = 10104
$ws.Range("M122").Value = -400
$ws.Range("N122").Value = -15004

$ws.Range("H132").Value = 9488430
$ws.Range("I132").Value = 9527365
$ws.Range("J132").Value = 9261309
$ws.Range("K132").Value = 28582095
$ws.Range("L132").Value = 27783927
$ws.Range("M132").Value = -28579565
$ws.Range("N132").Value = -27788987

$ws.Range("H136").Value = 2094432.2
$ws.Range("I136").Value = 1042577.06
$ws.Range("J136").Value = 7353708
$ws.Range("K136").Value = 3127731.18
$ws.Range("L136").Value = 22061124
$ws.Range("M136").Value = -3125181.18
$ws.Range("N136").Value = -22066224

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9081343
$ws.Range("I134").Value = 11629097
$ws.Range("J134").Value = 2234255
$ws.Range("K134").Value = 34887291
$ws.Range("L134").Value = 6702765
$ws.Range("M134").Value = -34884756
$ws.Range("N134").Value = -6707835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3250098
$ws.Range("I58").Value = 2385686.5
$ws.Range("J58").Value = 4546715
$ws.Range("K58").Value = 2385686.5
$ws.Range("L58").Value = 4546715
$ws.Range("M58").Value = -2385483.5
$ws.Range("N58").Value = -4547121

$ws.Range("H122").Value = 5484.5
$ws.Range("I122").Value = 9360.182000000001
$ws.Range("K122").Value = 28080.546
$ws.Range("M122").Value = -25630.546

$ws.Range("H132").Value = 1854416.6
$ws.Range("I132").Value = 3335228.5
$ws.Range("J132").Value = 3402
$ws.Range("K132").Value = 10005685.5
$ws.Range("L132").Value = 10206
$ws.Range("M132").Value = -10003155.5
$ws.Range("N132").Value = -15266

$ws.Range("H134").Value = 2676866.2
$ws.Range("I134").Value = 14104.5
$ws.Range("J134").Value = 5720022.5
$ws.Range("K134").Value = 42313.5
$ws.Range("L134").Value = 17160067.5
$ws.Range("M134").Value = -39778.5
$ws.Range("N134").Value = -17165137.5

$ws.Range("H136").Value = 3250098
$ws.Range("I136").Value = 2385686.5
$ws.Range("J136").Value = 4546715
$ws.Range("K136").Value = 7157059.5
$ws.Range("L136").Value = 13640145
$ws.Range("M136").Value = -7154509.5
$ws.Range("N136").Value = -13645245

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 1969.0769
$ws.Range("I125").Value = 1738
$ws.Range("J125").Value = 2167.1428
$ws.Range("K125").Value = 5214
$ws.Range("L125").Value = 6501.428400000001
$ws.Range("M125").Value = -294
$ws.Range("N125").Value = -16341.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2488426.2
$ws.Range("I70").Value = 1467589.1
$ws.Range("J70").Value = 4335655.5
$ws.Range("K70").Value = 1467589.1
$ws.Range("L70").Value = 4335655.5
$ws.Range("M70").Value = -1467319.1
$ws.Range("N70").Value = -4336195.5

$ws.Range("H73").Value = 2488426.2
$ws.Range("I73").Value = 1467589.1
$ws.Range("J73").Value = 4335655.5
$ws.Range("K73").Value = 1467589.1
$ws.Range("L73").Value = 4335655.5
$ws.Range("M73").Value = -1466653.1
$ws.Range("N73").Value = -4337527.5

$ws.Range("H80").Value = 16387.875
$ws.Range("I80").Value = 6362.5
$ws.Range("J80").Value = 26413.25
$ws.Range("K80").Value = 6362.5
$ws.Range("L80").Value = 26413.25
$ws.Range("M80").Value = -5364.5
$ws.Range("N80").Value = -28409.25

$ws.Range("H83").Value = 16387.875
$ws.Range("I83").Value = 6362.5
$ws.Range("J83").Value = 26413.25
$ws.Range("K83").Value = 31812.5
$ws.Range("L83").Value = 132066.25
$ws.Range("M83").Value = -26820.5
$ws.Range("N83").Value = -142050.25

$ws.Range("H122").Value = 9290677
$ws.Range("I122").Value = 39587.645
$ws.Range("K122").Value = 118762.935
$ws.Range("M122").Value = -116312.935

$ws.Range("H132").Value = 26043146
$ws.Range("I132").Value = 49524430
$ws.Range("J132").Value = 11367341
$ws.Range("K132").Value = 148573290
$ws.Range("L132").Value = 34102023
$ws.Range("M132").Value = -148570760
$ws.Range("N132").Value = -34107083

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1024.0358
$ws.Range("I7").Value = 808.05884
$ws.Range("J7").Value = 1357.8182
$ws.Range("K7").Value = 808.05884
$ws.Range("L7").Value = 1357.8182
$ws.Range("M7").Value = -696.05884
$ws.Range("N7").Value = -1581.8182

$ws.Range("H40").Value = 2085.9473
$ws.Range("I40").Value = 1186.9
$ws.Range("K40").Value = 1186.9
$ws.Range("M40").Value = -1050.9

$ws.Range("H82").Value = 3121.6333
$ws.Range("I82").Value = 1326.2
$ws.Range("J82").Value = 4917.067
$ws.Range("K82").Value = 1326.2
$ws.Range("L82").Value = 4917.067
$ws.Range("M82").Value = -965.2
$ws.Range("N82").Value = -5639.067

$ws.Range("H85").Value = 3121.6333
$ws.Range("I85").Value = 1326.2
$ws.Range("J85").Value = 4917.067
$ws.Range("K85").Value = 1326.2
$ws.Range("L85").Value = 4917.067
$ws.Range("M85").Value = -78.20000000000005
$ws.Range("N85").Value = -7413.067

$ws.Range("H122").Value = 6707717
$ws.Range("I122").Value = 927290.4399999999
$ws.Range("K122").Value = 2781871.32
$ws.Range("M122").Value = -2779421.32

$ws.Range("H126").Value = 1024.0358
$ws.Range("I126").Value = 808.05884
$ws.Range("J126").Value = 1357.8182
$ws.Range("K126").Value = 2424.17652
$ws.Range("L126").Value = 4073.4546
$ws.Range("M126").Value = 45.82348000000002
$ws.Range("N126").Value = -9013.454600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1245.1515
$ws.Range("I122").Value = 1164.4762
$ws.Range("K122").Value = 3493.4286
$ws.Range("M122").Value = -1043.4286

$ws.Range("H132").Value = 1545185.1
$ws.Range("I132").Value = 1152076.6
$ws.Range("J132").Value = 2527956.2
$ws.Range("K132").Value = 3456229.8
$ws.Range("L132").Value = 7583868.600000001
$ws.Range("M132").Value = -3453699.8
$ws.Range("N132").Value = -7588928.600000001
